$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.980.03"
$ws.Range("E2").Value = "  -2.01%  "
$ws.Range("D3").Value = "2.488.30"
$ws.Range("E3").Value = "  -3.70%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.22%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "302.41"
$ws.Range("E5").Value = "  -0.17%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "94.40"
$ws.Range("E6").Value = "  -2.72%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.577"
$ws.Range("E7").Value = "  +0.33%  "
$ws.Range("E8").Value = "  +0.33%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.529"
$ws.Range("E9").Value = "  -3.76%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.93"
$ws.Range("E10").Value = "  -2.41%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0801"
$ws.Range("E11").Value = "  -1.10%  "
$ws.Range("E12").Value = "  -2.26%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.35"
$ws.Range("E13").Value = "  -3.48%  "
$ws.Range("D14").Value = "2.879.28"
$ws.Range("E14").Value = "  -3.41%  "
$ws.Range("D15").Value = "2.505.30"
$ws.Range("E15").Value = "  -2.63%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.87"
$ws.Range("E16").Value = "  +3.39%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.849"
$ws.Range("E17").Value = "  -4.23%  "
$ws.Range("D18").Value = "42.146.03"
$ws.Range("E18").Value = "  -1.67%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.73"
$ws.Range("E19").Value = "  -1.14%  "
$ws.Range("D20").Value = "0.0₃0961"
$ws.Range("E20").Value = "  -2.92%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.34"
$ws.Range("E21").Value = "  -4.83%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "70.64"
$ws.Range("E22").Value = "  -1.81%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "246.52"
$ws.Range("E23").Value = "  -3.20%  "
$ws.Range("E24").Value = "  -2.65%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.99"
$ws.Range("E25").Value = "  -6.78%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.40"
$ws.Range("E26").Value = "  -7.86%  "
$ws.Range("E27").Value = "  -0.14%  "
$ws.Range("E28").Value = "  +8.47%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.03"
$ws.Range("E29").Value = "  -1.67%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "36.99"
$ws.Range("E30").Value = "  -5.12%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.84"
$ws.Range("E31").Value = "  -3.32%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "154.01"
$ws.Range("E32").Value = "  -0.92%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.25"
$ws.Range("E33").Value = "  -3.37%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.61"
$ws.Range("E34").Value = "  -5.44%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0774"
$ws.Range("E35").Value = "  -5.08%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.04"
$ws.Range("E36").Value = "  -6.37%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "18.21"
$ws.Range("E37").Value = "  -0.62%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.113"
$ws.Range("E38").Value = "  -1.33%  "
$ws.Range("E39").Value = "  -1.68%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "23.72"
$ws.Range("E40").Value = "  +1.45%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.80"
$ws.Range("E41").Value = "  -2.09%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.33"
$ws.Range("E42").Value = "  -2.58%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.999"
$ws.Range("E43").Value = "  -0.03%  "
$ws.Range("D44").Value = "2.037.71"
$ws.Range("E44").Value = "  -1.65%  "
$ws.Range("E45").Value = "  -4.90%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.96"
$ws.Range("E46").Value = "  -6.21%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.88"
$ws.Range("E47").Value = "  -4.15%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "82.56"
$ws.Range("E48").Value = "  -3.14%  "
$ws.Range("D49").Value = "2.736.70"
$ws.Range("E49").Value = "  -3.51%  "
$ws.Range("B50").Value = "ordi"
$ws.Range("C50").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "71.21"
$ws.Range("E50").Value = "  -6.58%  "
$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.186"
$ws.Range("E51").Value = "  -2.57%  "
